{"js": "// 1) \"Il se cultive en\" -> \"Il se cultive en \" (add a trailing space)\nconst cultiveHits = context.document.body.search(\"Il se cultive en\", { matchCase: true, matchWholeWord: false });\ncultiveHits.load(\"text\");\n\n// 2) \"que \" -> \"que si \" (only the FIRST occurrence in the document; a second,\n//    unrelated hit lives inside \"quelque \" further down and must stay untouched)\nconst queHits = context.document.body.search(\"que \", { matchCase: true, matchWholeWord: false });\nqueHits.load(\"text\");\n\n// 3) the literal tag text \"<sup>\" / \"</sup>\" -> \"<corr>\" / \"</corr>\", recoloured\n//    from #7F6000 to #A91111\nconst openTagHits = context.document.body.search(\"<sup>\", { matchCase: true, matchWholeWord: false });\nopenTagHits.load(\"text\");\n\nconst closeTagHits = context.document.body.search(\"</sup>\", { matchCase: true, matchWholeWord: false });\ncloseTagHits.load(\"text\");\n\nawait context.sync();\n\nif (cultiveHits.items.length > 0) {\n  cultiveHits.items[0].insertText(\"Il se cultive en \", \"Replace\");\n}\n\nif (queHits.items.length > 0) {\n  // items are in document order, so [0] is the real \"que \" (the other hit,\n  // inside \"quelque \", is items[1] and must be left alone)\n  queHits.items[0].insertText(\"que si \", \"Replace\");\n}\n\nif (openTagHits.items.length > 0) {\n  const r = openTagHits.items[0];\n  r.insertText(\"<corr>\", \"Replace\");\n  r.font.color = \"#A91111\";\n}\n\nif (closeTagHits.items.length > 0) {\n  const r = closeTagHits.items[0];\n  r.insertText(\"</corr>\", \"Replace\");\n  r.font.color = \"#A91111\";\n}\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# 1) \"Il se cultive en\" -> \"Il se cultive en \" (add a trailing space)\n$rng = $d.Content\nif ($rng.Find.Execute(\"Il se cultive en\")) {\n    $rng.Text = \"Il se cultive en \"\n}\n\n# 2) \"que \" -> \"que si \" -- search from the top of the document and take the\n#    FIRST hit only; a second, unrelated \"que \" lives inside \"quelque \" later\n#    on and must not be touched.\n$rng = $d.Content\nif ($rng.Find.Execute(\"que \")) {\n    $rng.Text = \"que si \"\n}\n\n# 3) literal tag text \"<sup>\" -> \"<corr>\", recoloured from #7F6000 to #A91111\n#    (Word/VBA colors are stored BGR: R + G*256 + B*65536 => 0xA9,0x11,0x11 = 1118633)\n$rng = $d.Content\nif ($rng.Find.Execute(\"<sup>\")) {\n    $rng.Text = \"<corr>\"\n    $rng.Font.Color = 1118633\n}\n\n# 4) literal tag text \"</sup>\" -> \"</corr>\", recoloured from #7F6000 to #A91111\n$rng = $d.Content\nif ($rng.Find.Execute(\"</sup>\")) {\n    $rng.Text = \"</corr>\"\n    $rng.Font.Color = 1118633\n}\n"}
